# feat: Control de cantidades especificas en carrito de ventas
#
# - Productos: el stock de "Americano" se actualiza (100 -> 50) y la fila se
#   reordena al final de la tabla (Latte y Capuccino quedan primero).
# - Ventas: se registran 3 nuevas ventas realizadas desde el carrito.
# - RegistroCaja: se reestructura con encabezados "Fecha/Hora" / "Operacion"
#   y una fila por cada venta registrada.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Hoja "Productos": reordenar Americano al final con el stock actualizado
# ---------------------------------------------------------------------
$wsProductos = $wb.Worksheets.Item("Productos")

# Quita la fila actual de "Americano" (antes en la fila 2)
$wsProductos.Range("A2:C2").ClearContents()

# La vuelve a agregar al final de la tabla con el stock actualizado
$wsProductos.Range("A5").Value = "Americano"
$wsProductos.Range("B5").Value = 4000
$wsProductos.Range("C5").Value = 50

# ---------------------------------------------------------------------
# Hoja "Ventas": nuevas ventas generadas desde el carrito
# ---------------------------------------------------------------------
$wsVentas = $wb.Worksheets.Item("Ventas")

$wsVentas.Range("A3").Value = "V-99E3C23B"
$wsVentas.Range("B3").Value = "2025-08-01 18:44:14"
$wsVentas.Range("C3").Value = 8000

$wsVentas.Range("A4").Value = "V-745239BB"
$wsVentas.Range("B4").Value = "2025-08-01 18:53:47"
$wsVentas.Range("C4").Value = 12000

$wsVentas.Range("A5").Value = "V-58AE91BC"
$wsVentas.Range("B5").Value = "2025-08-01 18:54:00"
$wsVentas.Range("C5").Value = 12000

# ---------------------------------------------------------------------
# Hoja "RegistroCaja": encabezados nuevos + una fila por venta
# ---------------------------------------------------------------------
$wsCaja = $wb.Worksheets.Item("RegistroCaja")

$wsCaja.Range("A1").Value = "Fecha/Hora"
$wsCaja.Range("B1").Value = "Operación"

$wsCaja.Range("A2").Value = "2025-08-01 18:44:14"
$wsCaja.Range("B2").Value = "VENTA - ID: V-99E3C23B | Monto: `$8000,00"

$wsCaja.Range("A3").Value = "2025-08-01 18:53:47"
$wsCaja.Range("B3").Value = "VENTA - ID: V-745239BB | Monto: `$12000,00"

$wsCaja.Range("A4").Value = "2025-08-01 18:54:00"
$wsCaja.Range("B4").Value = "VENTA - ID: V-58AE91BC | Monto: `$12000,00"
